# Updated cryptos list - applies new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.121.53'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.788.01'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.78'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.546'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.96'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0690'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('D12').Value = '2.045.92'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.19'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').Value = '1.787.48'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '34.064.74'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.17'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.06'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.38'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.15'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.31'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.65'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.60'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('D35').Value = '1.457.31'
$ws.Range('E35').Value = '  +4.36%  '
$ws.Range('E36').Value = '  +9.17%  '
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('E38').Value = '  +2.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.03'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.34'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.81%  '
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.51'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.06'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0505'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.79'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('D50').Value = '1.946.94'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('E51').Value = '  +0.12%  '
